$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2826 (EMA6/EMA12/EMA24 recomputed for the new window) ---
$ws.Range("G2826").Value = 77.56215221316754
$ws.Range("H2826").Value = 77.90298529819698
$ws.Range("I2826").Value = 77.93080223674818

# --- Append 50 new rows (2827-2876) of windowed OHLC + EMA + trades data ---
$rows = @(
    @(2827,113,44620,75.83000183105469,77.86000061035156,74.77999877929688,77.83000183105469,77.56215221316754,77.90298529819698,77.93080223674818,114),
    @(2828,114,44627,77.01000213623047,79.62999725341797,76.08000183105469,78.26000213623047,77.76153790547123,77.95791096558675,77.95713822870677,114),
    @(2829,115,44634,77.90000152587891,79.73999786376953,77.30999755859375,79.11000061035156,78.14681296400846,78.13515552631979,78.04936721923835,114),
    @(2830,116,44641,79.66999816894531,81.37999725341797,78.76999664306641,81.33999633789062,79.05915107083193,78.62820795886915,78.31261754873053,114),
    @(2831,117,44648,81.37000274658203,83.58000183105469,80.73000335693359,83.51999664306641,80.33367837718464,79.38079083336103,78.7292078762774,114),
    @(2832,118,44655,83.5,87.83999633789062,82.73000335693359,87.68000030517578,82.43262749946783,80.6575922905633,79.44527127058927,114),
    @(2833,119,44662,88.31999969482422,89.48000335693359,84.83000183105469,86.91000366210938,83.71187783165112,81.61950173233961,80.04244986191088,114),
    @(2834,120,44669,86.41999816894531,87.44999694824219,84.5,84.58999633789062,83.96276883343384,82.07650090242439,80.40625357998925,114),
    @(2835,121,44676,84.58999633789062,90.01000213623047,83.52999877929688,88.69000244140625,85.31340700714024,83.09396267765237,81.06895348890261,114),
    @(2836,122,44683,88.72000122070312,88.91999816894531,86.36000061035156,88.38999938964844,86.19243340214258,83.908737556421,81.65463716096228,114),
    @(2837,123,44690,87.5,91.05000305175781,86.91000366210938,90.41000366210938,87.3974534764188,84.90893234191152,82.35506648105405,114),
    @(2838,124,44697,90.79000091552734,94.55999755859375,90.68000030517578,93.55000305175781,89.15532478365851,86.23832783573403,83.25066140671035,114),
    @(2839,125,44704,93.58999633789062,94.91999816894531,91,93.08000183105469,90.27666108291456,87.29089306578337,84.0370086406579,114),
    @(2840,126,44711,91.91000366210938,92.62000274658203,88.73999786376953,89.91000366210938,90.17190181982737,87.69383315752583,84.50684824237402,114),
    @(2841,127,44718,90.16999816894531,90.81999969482422,86.70999908447266,87.18000030517578,89.3170728156412,87.61478194947198,84.72070040739816,114),
    @(2842,128,44725,85.70999908447266,86.51000213623047,83.05000305175781,84.62000274658203,87.97505279591,87.1540466874889,84.71264459453286,114),
    @(2843,129,44732,86.05000305175781,93.79000091552734,85.26999664306641,93.12999725341797,89.44789406948371,88.07342369763184,85.38603280724367,114),
    @(2844,130,44739,93.19999694824219,95.72000122070312,90.38999938964844,92.41999816894531,90.29706666932988,88.74212746244929,85.9487500361798,114),
    @(2845,131,44746,92.44000244140625,94.51000213623047,90.52999877929688,92.77999877929688,91.00647584360617,89.36333843427199,86.49524993562918,114),
    @(2846,132,44753,92.87999725341797,95.34999847412109,92.05999755859375,94.95999908447266,92.13605391242517,90.22436314968748,87.17242986753665,114),
    @(2847,133,44760,94.90000152587891,95.09999847412109,88.91000366210938,90.11000061035156,91.55718154040414,90.20676891286656,87.40743552696185,114),
    @(2848,134,44767,90.19999694824219,91.94999694824219,87.41999816894531,89.33999633789062,90.92370005397171,90.0734192859472,87.56204039183615,114),
    @(2849,135,44774,89.33999633789062,90.01999664306641,86.01999664306641,87.41000366210938,89.9197867991539,89.663663036126,87.54987745345801,114),
    @(2850,136,44781,87.30000305175781,91.02999877929688,86.80999755859375,91.01999664306641,90.23413246884319,89.87232974488606,87.82748698862667,114),
    @(2851,0,"Month",90.83000183105469,93.01999664306641,85.06999969482422,87.33999633789062,77.76153790547123,77.95791096558675,77.95713822870677,114),
    @(2852,114,44627,77.01000213623047,79.62999725341797,76.08000183105469,78.26000213623047,77.76153790547123,77.95791096558675,77.95713822870677,115),
    @(2853,115,44634,77.90000152587891,79.73999786376953,77.30999755859375,79.11000061035156,78.14681296400846,78.13515552631979,78.04936721923835,115),
    @(2854,116,44641,79.66999816894531,81.37999725341797,78.76999664306641,81.33999633789062,79.05915107083193,78.62820795886915,78.31261754873053,115),
    @(2855,117,44648,81.37000274658203,83.58000183105469,80.73000335693359,83.51999664306641,80.33367837718464,79.38079083336103,78.7292078762774,115),
    @(2856,118,44655,83.5,87.83999633789062,82.73000335693359,87.68000030517578,82.43262749946783,80.6575922905633,79.44527127058927,115),
    @(2857,119,44662,88.31999969482422,89.48000335693359,84.83000183105469,86.91000366210938,83.71187783165112,81.61950173233961,80.04244986191088,115),
    @(2858,120,44669,86.41999816894531,87.44999694824219,84.5,84.58999633789062,83.96276883343384,82.07650090242439,80.40625357998925,115),
    @(2859,121,44676,84.58999633789062,90.01000213623047,83.52999877929688,88.69000244140625,85.31340700714024,83.09396267765237,81.06895348890261,115),
    @(2860,122,44683,88.72000122070312,88.91999816894531,86.36000061035156,88.38999938964844,86.19243340214258,83.908737556421,81.65463716096228,115),
    @(2861,123,44690,87.5,91.05000305175781,86.91000366210938,90.41000366210938,87.3974534764188,84.90893234191152,82.35506648105405,115),
    @(2862,124,44697,90.79000091552734,94.55999755859375,90.68000030517578,93.55000305175781,89.15532478365851,86.23832783573403,83.25066140671035,115),
    @(2863,125,44704,93.58999633789062,94.91999816894531,91,93.08000183105469,90.27666108291456,87.29089306578337,84.0370086406579,115),
    @(2864,126,44711,91.91000366210938,92.62000274658203,88.73999786376953,89.91000366210938,90.17190181982737,87.69383315752583,84.50684824237402,115),
    @(2865,127,44718,90.16999816894531,90.81999969482422,86.70999908447266,87.18000030517578,89.3170728156412,87.61478194947198,84.72070040739816,115),
    @(2866,128,44725,85.70999908447266,86.51000213623047,83.05000305175781,84.62000274658203,87.97505279591,87.1540466874889,84.71264459453286,115),
    @(2867,129,44732,86.05000305175781,93.79000091552734,85.26999664306641,93.12999725341797,89.44789406948371,88.07342369763184,85.38603280724367,115),
    @(2868,130,44739,93.19999694824219,95.72000122070312,90.38999938964844,92.41999816894531,90.29706666932988,88.74212746244929,85.9487500361798,115),
    @(2869,131,44746,92.44000244140625,94.51000213623047,90.52999877929688,92.77999877929688,91.00647584360617,89.36333843427199,86.49524993562918,115),
    @(2870,132,44753,92.87999725341797,95.34999847412109,92.05999755859375,94.95999908447266,92.13605391242517,90.22436314968748,87.17242986753665,115),
    @(2871,133,44760,94.90000152587891,95.09999847412109,88.91000366210938,90.11000061035156,91.55718154040414,90.20676891286656,87.40743552696185,115),
    @(2872,134,44767,90.19999694824219,91.94999694824219,87.41999816894531,89.33999633789062,90.92370005397171,90.0734192859472,87.56204039183615,115),
    @(2873,135,44774,89.33999633789062,90.01999664306641,86.01999664306641,87.41000366210938,89.9197867991539,89.663663036126,87.54987745345801,115),
    @(2874,136,44781,87.30000305175781,91.02999877929688,86.80999755859375,91.01999664306641,90.23413246884319,89.87232974488606,87.82748698862667,115),
    @(2875,137,44788,90.83000183105469,93.01999664306641,89.87000274658203,92.08000183105469,90.76152371518933,90.21197160429662,88.16768817602092,115),
    @(2876,0,"Month",91.80999755859375,92.22000122070312,85.06999969482422,87.33999633789062,90.76152371518933,90.21197160429662,88.16768817602092,115)
)

foreach ($row in $rows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]

    # Column A: bold, centered, bordered header-like style (matches the rest of the sheet's "A" column)
    $ws.Cells.Item($r, 1).Font.Bold = $true
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 1).VerticalAlignment = -4160
    $ws.Cells.Item($r, 1).Borders.LineStyle = 1
    $ws.Cells.Item($r, 1).Borders.Weight = 2

    # Column B: date-time format, except on "Month" summary rows where B holds text
    if ($row[2] -ne "Month") {
        $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
}

Write-Output "applied edit: added 50 rows (2827-2876) + updated row 2826"
